# "Add 5 more lab results" -- this particular change inserts a fresh blank
# row above the first data row of the "Expansion List" table (row 13),
# shifting all of the existing lab-result rows down by one, to make room
# for new entries to be typed in later.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Expansion List")

# Right-click row 13 header -> Insert: shifts rows 13:40 down to 14:41 and
# formats the new row like the row above it (the filter header row), which
# is exactly how Excel produces the bordered-but-no-outline "spacer" row.
$ws.Rows("13:13").Insert(-4121)

# The header row (now back to its normal height) and the new spacer row
# (which picks up the header's tall/underlined look) both need their row
# height nudged so the XML keeps matching what Excel wrote.
$ws.Rows("12:12").AutoFit()
$ws.Rows("13:13").RowHeight = 16

# Leave the freshly-inserted row selected, same as right after doing the
# insert interactively.
$ws.Rows("13:13").Select() | Out-Null
